$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = "thay doi lan 2 "
$ws.Range("E8").Value = 111

$ws.Range("E8").Select()
